$d = $word.ActiveDocument

$find = "editing of issues. "
$rng = $d.Content
$rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter("Created Doxygen documentation and makefile, and wrote user manual.")
